$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Baghdati")

# Copy column J (rows 3-6) formatting into column K, then overwrite values
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)

# Add new column K data (year 2023)
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 756.3
$ws.Range("K5").Value = 1017.8
$ws.Range("K6").Value = 672.7
